# Applies the change described by the diff:
#  - Updates the date/time serial values in column A (rows 2-27) of the
#    active sheet to a new set of timestamps (2024-05-19 22:00 through
#    2024-05-20 22:00, with tiny float drift matching the source data).
#  - Updates the sheet view: topLeftCell="A4", selection activeCell="A4"
#    with sqref "A4:A27".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    45431.916666666664,
    45431.958333333336,
    45432,
    45432.041666666664,
    45432.08333321759,
    45432.124999826388,
    45432.166666435187,
    45432.208333043978,
    45432.249999652777,
    45432.291666261575,
    45432.333332870374,
    45432.374999479165,
    45432.416666087964,
    45432.458332696762,
    45432.499999305554,
    45432.541665914352,
    45432.583332523151,
    45432.624999131942,
    45432.66666574074,
    45432.708332349539,
    45432.74999895833,
    45432.791665567129,
    45432.833332175927,
    45432.874998784719,
    45432.916665393517,
    45432.958332002316
)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 1).Value = $v
    $row = $row + 1
}

# Update the sheet view: scroll so row 4 is the top-left visible cell, and
# select A4:A27 with A4 as the active cell (selecting the range sets
# A4 as the active cell automatically, since it's the first cell of the
# selection).
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A4:A27").Select()
